$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova coluna B (dados de 2020) e desloca as colunas
# existentes (2021-2025) uma posicao para a direita.
$ws.Columns("B:B").Insert()

# Cabecalho (anos)
$ws.Range("B1").Value = 2020
$ws.Range("C1").Value = 2021
$ws.Range("D1").Value = 2022
$ws.Range("E1").Value = 2023
$ws.Range("F1").Value = 2024
$ws.Range("G1").Value = 2025

# Linha 2 (mes 1)
$ws.Range("C2").Value = 11682.62
$ws.Range("D2").Value = 43029.98
$ws.Range("E2").Value = 52760.71
$ws.Range("F2").Value = 90737.31
$ws.Range("G2").Value = 106156.53

# Linha 3 (mes 2)
$ws.Range("C3").Value = 8654.09
$ws.Range("D3").Value = 42931.49
$ws.Range("E3").Value = 55235.92
$ws.Range("F3").Value = 93714.89
$ws.Range("G3").Value = 110632.57

# Linha 4 (mes 3)
$ws.Range("C4").Value = 13401.71
$ws.Range("D4").Value = 46060.74
$ws.Range("E4").Value = 64193.54
$ws.Range("F4").Value = 84818.77
$ws.Range("G4").Value = 83990.55

# Linha 5 (mes 4)
$ws.Range("C5").Value = 16166.33
$ws.Range("D5").Value = 42566.03
$ws.Range("E5").Value = 60035.67
$ws.Range("F5").Value = 101395.59
$ws.Range("G5").Value = 95827.67

# Linha 6 (mes 5)
$ws.Range("C6").Value = 13436.58
$ws.Range("D6").Value = 36440.33
$ws.Range("E6").Value = 68445.73
$ws.Range("F6").Value = 92789.5

# Linha 7 (mes 6)
$ws.Range("C7").Value = 18576.1
$ws.Range("D7").Value = 41979.5
$ws.Range("E7").Value = 68280.27
$ws.Range("F7").Value = 101320.28

# Linha 8 (mes 7)
$ws.Range("C8").Value = 15605.42
$ws.Range("D8").Value = 36927.96
$ws.Range("E8").Value = 63709.16
$ws.Range("F8").Value = 103779.05

# Linha 9 (mes 8)
$ws.Range("C9").Value = 30073.31
$ws.Range("D9").Value = 39097.26
$ws.Range("E9").Value = 77793.75
$ws.Range("F9").Value = 102460.71

# Linha 10 (mes 9)
$ws.Range("B10").Value = 5513.09
$ws.Range("C10").Value = 32873.05
$ws.Range("D10").Value = 50202.25
$ws.Range("E10").Value = 77438.9
$ws.Range("F10").Value = 108365.58

# Linha 11 (mes 10)
$ws.Range("B11").Value = 5255.93
$ws.Range("C11").Value = 42033.55
$ws.Range("D11").Value = 38621.84
$ws.Range("E11").Value = 85404.96
$ws.Range("F11").Value = 124976.98

# Linha 12 (mes 11)
$ws.Range("B12").Value = 8403.630000000001
$ws.Range("C12").Value = 42634.32
$ws.Range("D12").Value = 41810.82
$ws.Range("E12").Value = 64809.84
$ws.Range("F12").Value = 80778.8

# Linha 13 (mes 12)
$ws.Range("B13").Value = 4693.76
$ws.Range("C13").Value = 34781.06
$ws.Range("D13").Value = 44681.08
$ws.Range("E13").Value = 64356.44
$ws.Range("F13").Value = 111925.01
